# "Refined metadata to be additional tab"
#
# 1) Refresh the panel-query timestamps on the "data" sheet (column F,
#    rows 2-7) to reflect the newer query run.
# 2) Pull the per-panel metadata (name/id/version/query time/request
#    URL) out into its own "metadata" worksheet, placed right after
#    "data".

$wb   = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1) Refresh the query timestamps on the existing "data" sheet
# ---------------------------------------------------------------------
$data.Range("F2").Value = "2021-10-05 14:22:12.600132"
$data.Range("F3").Value = "2021-10-05 14:22:12.600140"
$data.Range("F4").Value = "2021-10-05 14:22:12.600143"
$data.Range("F5").Value = "2021-10-05 14:22:12.600146"
$data.Range("F6").Value = "2021-10-05 14:22:12.600149"
$data.Range("F7").Value = "2021-10-05 14:22:12.600151"

# ---------------------------------------------------------------------
# 2) Add the new "metadata" tab, placed right after "data"
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Give the header cells the same bold/boxed look as the "data" sheet's
# header row by copying its formatting over (keeps the shared style
# table tidy instead of minting new font/xf entries).
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("B1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

# Data row - column A uses the same index-column style as "data"
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "PHACE(S) syndrome"
$meta.Range("C2").Value = 94

# "1.2" must stay a text value (not be coerced into the number 1.2):
# force text entry, then drop back to the default "Normal" style so it
# matches the un-styled look of the rest of the data row.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.2"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2020-05-07T14:26:47.280928Z"
$meta.Range("F2").Value = "2021-10-05 14:22:12.596943"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/94/?format=json"
